$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "299.47"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "2.25%"
$ws.Cells.Item(2,5).Style = "Normal"

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "42.18"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "4.70%"
$ws.Cells.Item(3,5).Style = "Normal"

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "5.013"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "0.22%"
$ws.Cells.Item(4,5).Style = "Normal"

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.07561"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "3.20%"
$ws.Cells.Item(5,5).Style = "Normal"

$ws.Cells.Item(6,2).Value = "GateToken"
$ws.Cells.Item(6,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "4.381"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "2.03%"
$ws.Cells.Item(6,5).Style = "Normal"

$ws.Cells.Item(7,2).Value = "FTXToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "1.599"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "2.42%"
$ws.Cells.Item(7,5).Style = "Normal"

$ws.Cells.Item(8,2).Value = "MXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.9385"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "1.79%"
$ws.Cells.Item(8,5).Style = "Normal"

$ws.Cells.Item(9,2).Value = "BTSEToken"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "2.384"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "0.14%"
$ws.Cells.Item(9,5).Style = "Normal"

$ws.Cells.Item(10,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.1189"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "0.20%"
$ws.Cells.Item(10,5).Style = "Normal"

$ws.Cells.Item(11,2).Value = "WazirX"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.1836"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "1.16%"
$ws.Cells.Item(11,5).Style = "Normal"

$ws.Cells.Item(12,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.09043"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "2.66%"
$ws.Cells.Item(12,5).Style = "Normal"

$ws.Cells.Item(13,2).Value = "BitrueCoin"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.04167"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "-5.16%"
$ws.Cells.Item(13,5).Style = "Normal"

$ws.Cells.Item(14,2).Value = "BitMartToken"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.1048"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "-0.54%"
$ws.Cells.Item(14,5).Style = "Normal"

$ws.Cells.Item(15,2).Value = "BitForexToken"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.001282"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "0.38%"
$ws.Cells.Item(15,5).Style = "Normal"

$ws.Cells.Item(16,2).Value = "TigerCash"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.005802"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "0.57%"
$ws.Cells.Item(16,5).Style = "Normal"

$ws.Cells.Item(17,2).Value = "LEO"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "3.344"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "0.06%"
$ws.Cells.Item(17,5).Style = "Normal"

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "0.3334"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "0.19%"
$ws.Cells.Item(18,5).Style = "Normal"

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "8.384"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "6.58%"
$ws.Cells.Item(19,5).Style = "Normal"

$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "1.32%"
$ws.Cells.Item(20,5).Style = "Normal"

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.3299"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "17.71%"
$ws.Cells.Item(21,5).Style = "Normal"

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "0.04101"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "4.34%"
$ws.Cells.Item(22,5).Style = "Normal"

$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = "0.30%"
$ws.Cells.Item(23,5).Style = "Normal"

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "0.003902"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = "2.41%"
$ws.Cells.Item(24,5).Style = "Normal"

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.0001269"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = "1.46%"
$ws.Cells.Item(25,5).Style = "Normal"

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.02410"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = "3.14%"
$ws.Cells.Item(38,5).Style = "Normal"

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.05231"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = "3.10%"
$ws.Cells.Item(39,5).Style = "Normal"

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.006649"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = "10.56%"
$ws.Cells.Item(40,5).Style = "Normal"

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.007707"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = "-1.32%"
$ws.Cells.Item(41,5).Style = "Normal"

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.1327"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = "2.96%"
$ws.Cells.Item(42,5).Style = "Normal"

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.007392"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = "0.10%"
$ws.Cells.Item(43,5).Style = "Normal"

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.007812"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = "-2.95%"
$ws.Cells.Item(44,5).Style = "Normal"

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.3000"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = "3.20%"
$ws.Cells.Item(45,5).Style = "Normal"

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.00006252"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = "0.65%"
$ws.Cells.Item(46,5).Style = "Normal"

$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = "-0.18%"
$ws.Cells.Item(47,5).Style = "Normal"

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.04574"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = "-4.13%"
$ws.Cells.Item(48,5).Style = "Normal"

$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = "-0.09%"
$ws.Cells.Item(49,5).Style = "Normal"

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "0.00002098"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = "-0.18%"
$ws.Cells.Item(50,5).Style = "Normal"

$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = "-0.18%"
$ws.Cells.Item(51,5).Style = "Normal"
